# Poland Ekstraklasa - base update (19-04-2024 23:27)
# Adds two new matches (ids 251 and 252) at the top of the 2024-04-xx block,
# removes an outdated/duplicate match (B=6775594), and appends the match
# that used to close the block (B=6885526) as a brand-new last row.
#
# Net effect vs. the original sheet:
#   - rows 253..259 (B: 6775597,6774472,6775594,6850054,6830603,6775596,6885526)
#   become
#   - rows 253..260 (B: 6775595,6775598,6775597,6774472,6850054,6830603,6775596,6885526)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

function Get-RowVals($r) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Range("$c$r").Value2
    }
    return $vals
}

function Set-RowVals($r, $vals) {
    foreach ($c in $cols) {
        $v = $vals[$c]
        if ($null -ne $v) {
            $ws.Range("$c$r").Value = $v
        }
    }
}

# --- 1. Capture the existing rows we still need, before anything is overwritten ---
$old253 = Get-RowVals 253   # -> becomes new row 255 (A: 251 -> 253)
$old254 = Get-RowVals 254   # -> becomes new row 256 (A: 252 -> 254)
# old 255 (B=6775594) is dropped entirely
$old256 = Get-RowVals 256   # -> becomes new row 257 (A: 254 -> 255)
$old257 = Get-RowVals 257   # -> becomes new row 258 (A: 255 -> 256)
$old258 = Get-RowVals 258   # -> becomes new row 259 (A: 256 -> 257)
$old259 = Get-RowVals 259   # -> becomes new row 260 (A: 257 -> 258), brand new row

# Row 260 does not exist yet, so it has no cell formatting at all. Clone the
# "id" (bold + border) and "Date" (custom date format) styles from the row
# directly above it before filling in its values, matching the look of
# every other data row.
$ws.Range("A259").Copy()
$ws.Range("A260").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E259").Copy()
$ws.Range("E260").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- 2. Write the shifted rows (255..260), fixing up only the id column A ---
$old253["A"] = 253
Set-RowVals 255 $old253

$old254["A"] = 254
Set-RowVals 256 $old254

$old256["A"] = 255
Set-RowVals 257 $old256

$old257["A"] = 256
Set-RowVals 258 $old257

$old258["A"] = 257
Set-RowVals 259 $old258

$old259["A"] = 258
Set-RowVals 260 $old259

# --- 3. Write the two brand-new rows (253, 254) with their own fresh data ---
$new253 = @{
    A = 251; B = 6775595; C = "Poland Ekstraklasa"; D = "Poland Ekstraklasa";
    E = 45401.54166666666; F = "Korona Kielce"; G = "Radomiak Radom";
    H = 4; I = 0; J = "H";
    K = 2.375; L = 3.2; M = 3.1; N = 2.375; O = 3.2; P = 3.1;
    Q = -0.25; R = 2.025; S = 1.825; T = 2.25; U = 2.05; V = 1.8;
    W = 1.375; X = -1; Y = -1; Z = 1.025; AA = -1; AB = 1.05; AC = -1
}
Set-RowVals 253 $new253

$new254 = @{
    A = 252; B = 6775598; C = "Poland Ekstraklasa"; D = "Poland Ekstraklasa";
    E = 45401.64583333334; F = "Rakow Czestochowa"; G = "Gornik Zabrze";
    H = 0; I = 1; J = "A";
    K = 1.571; L = 4; M = 5.75; N = 1.533; O = 4.2; P = 6.5;
    Q = -1; R = 1.975; S = 1.875; T = 2.5; U = 1.85; V = 2;
    W = -1; X = -1; Y = 5.5; Z = -1; AA = 0.875; AB = -1; AC = 1
}
Set-RowVals 254 $new254
